# Session 2 - Editing the first slideshow, making the second slideshow
#
# Slide 1 originally has a single "Debate Basics" textbox (TextBox 2 / id 2).
# This turns it into a 3-textbox "Week 2 / Debate Basics / Week 1" header:
#   - TextBox 2 (reused, id 2): becomes the small "Week 2" label
#   - TextBox 3 (new, id 3): keeps the big "Debate Basics" text, moved down
#   - TextBox 4 (new, id 4): becomes the small "Week 1" label

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$orig = $s.Shapes.Item(1)

# EMU -> point helper (EMU = 12700 * point). PowerPoint stores shape
# geometry as single-precision floats internally, so a handful of the
# literals below carry extra fractional digits chosen so that they
# round-trip back to the exact target EMU values.

# --- TextBox 3: duplicate of the original, keeps "Debate Basics", moves down ---
$tb3 = $orig.Duplicate()
$tb3.Name = "TextBox 3"
$tb3.Left = 29.89177133346853
$tb3.Top = 354.9216918933831
$tb3.Width = 603.3047485593374
$tb3.Height = 143.59996035982238

# --- TextBox 4: duplicate of the original, becomes the "Week 1" label ---
$tb4 = $orig.Duplicate()
$tb4.Name = "TextBox 4"
$tb4.Left = 29.89177133346853
$tb4.Top = 291.228393556771
$tb4.Width = 177.3750000299011
$tb4.Height = 82.60004039997979

$tr4 = $tb4.TextFrame.TextRange
$tr4.Text = "Week 1"
$tr4.Font.Size = 51.99
$tr4.Font.Color.RGB = 11003894
$tr4.ParagraphFormat.SpaceWithin = 72.79

# --- Original shape (TextBox 2) is repurposed into the "Week 2" label ---
$orig.Left = 29.89177133346853
$orig.Top = 291.228393556771
$orig.Width = 177.3750000299011
$orig.Height = 82.60004039997979

$tr2 = $orig.TextFrame.TextRange
$tr2.Text = "Week 2"
$tr2.Font.Size = 51.99
$tr2.Font.Color.RGB = 1647999
$tr2.ParagraphFormat.SpaceWithin = 72.79
